# "update data with resort sheetname"
# Re-sort the worksheet tabs: move "总计" (the summary sheet) so it becomes
# the first sheet in the workbook, ahead of "2022-Q2". This only changes the
# sheet order - the data on each sheet is untouched - and keeps "2022-Q2"
# (which was the active sheet beforehand) selected afterwards, since moving
# a different tab shouldn't change which sheet the user had open.
$wb = $excel.ActiveWorkbook

$summarySheet = $wb.Worksheets.Item("总计")
$summarySheet.Move($wb.Worksheets.Item(1))

# Look the sheet back up by name (a reference captured before the Move
# above goes stale) so re-activating it really restores "2022-Q2" as the
# selected tab - moving a *different* sheet's tab shouldn't steal the
# active-sheet selection.
$wb.Worksheets.Item("2022-Q2").Activate()
